$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Entity Relationship Diagram -> add version v0.1
$ws.Range("B4").Value = "v0.1"

# Row 5: Data Flow Diagram -> add version v0.1
$ws.Range("B5").Value = "v0.1"

# Row 18: Class Analysis Diagram -> add version v1.0
$ws.Range("B18").Value = "v1.0"

# Row 19: IDEF0 -> add version v0.1
$ws.Range("B19").Value = "v0.1"

# Row 20: new entry - UML Activity Diagram, v1.0
$ws.Range("A20").Value = "UML Activity Diagram"
$ws.Range("B20").Value = "v1.0"

# Row 21: new entry - Техническое задание (ГОСТ 34.602-2020), рыба
$ws.Range("A21").Value = "Техническое задание (ГОСТ 34.602-2020)"
$ws.Range("B21").Value = "рыба"

# Match the author's final cursor position noted in the saved file
$ws.Range("B8").Select() | Out-Null
